$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.435.43'
$ws.Range('E2').Value = '  +1.83%  '
$ws.Range('D3').Value = '1.842.45'
$ws.Range('E3').Value = '  +1.68%  '
$ws.Range('E4').Value = '  +0.67%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '234.76'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +4.26%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.621'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.90%  '
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '43.37'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +10.02%  '
$ws.Range('E9').Value = '  +6.77%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0695'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.37%  '
$ws.Range('E11').Value = '  +0.97%  '
$ws.Range('D12').Value = '2.110.02'
$ws.Range('E12').Value = '  +1.62%  '
$ws.Range('D13').Value = '1.855.98'
$ws.Range('E13').Value = '  +2.48%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '11.30'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.90%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.76'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +8.04%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.671'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +5.41%  '
$ws.Range('D17').Value = '35.470.93'
$ws.Range('E17').Value = '  +1.99%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '70.57'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +3.43%  '
$ws.Range('D19').Value = '0.0₃0799'
$ws.Range('E19').Value = '  +3.73%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '241.91'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.98'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +7.50%  '
$ws.Range('E22').Value = '  +12.84%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.01'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('E24').Value = '  +3.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '171.28'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.90'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.21%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '17.69'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.12%  '
$ws.Range('E28').Value = '  +0.88%  '
$ws.Range('E29').Value = '  +33.39%  '
$ws.Range('E30').Value = '  +0.56%  '
$ws.Range('D31').Value = '3.342.34'
$ws.Range('E31').Value = '  +37.56%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0563'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +9.29%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.10'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +6.07%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.96'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.96%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.81'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '94.49'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +14.18%  '
$ws.Range('E37').Value = '  +6.82%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.13'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +6.40%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0195'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.77%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '15.38'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +5.33%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '1.322.72'
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.01'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +6.21%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.28'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.05%  '
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.45'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.11%  '
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '6.27'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +9.30%  '
$ws.Range('E48').Value = '  -0.93%  '
$ws.Range('D49').Value = '2.019.84'
$ws.Range('E49').Value = '  +2.18%  '
$ws.Range('E50').Value = '  +0.69%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '101.80'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.43%  '
